$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set ExisUnits (column F) to 0 for rows whose value changed in the source edit
$rowsToZeroF = @(8, 10, 12, 14, 15, 16)
foreach ($r in $rowsToZeroF) {
    $ws.Range("F$r").Value = 0
}

# Set MaxInvest (column I) to 200 for rows 8 through 18 (MaxlineLoad 100%)
for ($r = 8; $r -le 18; $r++) {
    $ws.Range("I$r").Value = 200
}

# Update the active selection to match the saved view state (I23:I24, active I23)
$ws.Range("I23:I24").Select()
